$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = '@'
    $rng.Value2 = $text
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '25.699.55'
Set-TextValue 'E2' '  -3.55%  '

Set-TextValue 'D3' '1.741.83'
Set-TextValue 'E3' '  -5.59%  '

Set-TextValue 'D4' '1.001'
Set-TextValue 'E4' '  +0.14%  '

Set-TextValue 'D5' '237.89'
Set-TextValue 'E5' '  -8.65%  '

Set-TextValue 'D6' '1.002'
Set-TextValue 'E6' '  +0.15%  '

Set-TextValue 'D7' '0.4909'
Set-TextValue 'E7' '  -6.99%  '

Set-TextValue 'D8' '41.87'
Set-TextValue 'E8' '  -7.18%  '

Set-TextValue 'D9' '0.2405'
Set-TextValue 'E9' '  -23.87%  '

Set-TextValue 'D10' '0.05987'
Set-TextValue 'E10' '  -11.96%  '

Set-TextValue 'D11' '1.739.91'
Set-TextValue 'E11' '  -5.69%  '

Set-TextValue 'D12' '0.06815'
Set-TextValue 'E12' '  -12.22%  '

Set-TextValue 'D13' '14.68'
Set-TextValue 'E13' '  -21.84%  '

Set-TextValue 'D14' '4.428'
Set-TextValue 'E14' '  -11.69%  '

Set-TextValue 'D15' '76.53'
Set-TextValue 'E15' '  -13.09%  '

Set-TextValue 'D16' '0.5774'
Set-TextValue 'E16' '  -26.38%  '

Set-TextValue 'D17' '1.001'
Set-TextValue 'E17' '  +0.15%  '

Set-TextValue 'E18' '  +0.10%  '

Set-TextValue 'D19' '25.729.90'
Set-TextValue 'E19' '  -3.55%  '

Set-TextValue 'D20' '11.45'
Set-TextValue 'E20' '  -17.39%  '

Set-TextValue 'D21' '0.000006393'
Set-TextValue 'E21' '  -19.27%  '

Set-TextValue 'D22' '1.960.01'
Set-TextValue 'E22' '  -5.65%  '

Set-TextValue 'D23' '3.944'
Set-TextValue 'E23' '  -14.43%  '

Set-TextValue 'D24' '5.079'
Set-TextValue 'E24' '  -15.08%  '

Set-TextValue 'D25' '7.830'
Set-TextValue 'E25' '  -15.89%  '

Set-TextValue 'D26' '136.28'
Set-TextValue 'E26' '  -4.38%  '

Set-TextValue 'B27' 'Toncoin'
Set-TextValue 'C27' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D27' '1.470'
Set-TextValue 'E27' '  -12.71%  '

Set-TextValue 'B28' 'LidoDAOToken'
Set-TextValue 'C28' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D28' '1.843'
Set-TextValue 'E28' '  -16.88%  '

Set-TextValue 'D29' '14.48'
Set-TextValue 'E29' '  -14.93%  '

Set-TextValue 'D30' '100.01'
Set-TextValue 'E30' '  -9.84%  '

Set-TextValue 'B31' 'InternetComputer(DFINITY)'
Set-TextValue 'C31' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D31' '3.748'
Set-TextValue 'E31' '  -10.80%  '

Set-TextValue 'B32' 'Stellar'
Set-TextValue 'C32' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D32' '0.08096'
Set-TextValue 'E32' '  -7.25%  '

Set-TextValue 'D33' '3.353'
Set-TextValue 'E33' '  -18.00%  '

Set-TextValue 'D34' '0.04376'
Set-TextValue 'E34' '  -10.49%  '

Set-TextValue 'E35' '  +0.14%  '

Set-TextValue 'D36' '2.697'
Set-TextValue 'E36' '  -5.60%  '

Set-TextValue 'D37' '1.013'
Set-TextValue 'E37' '  -11.31%  '

Set-TextValue 'D38' '0.5997'
Set-TextValue 'E38' '  -18.12%  '

Set-TextValue 'D39' '2.719'
Set-TextValue 'E39' '  -12.38%  '

Set-TextValue 'D40' '2.074'
Set-TextValue 'E40' '  -9.07%  '

Set-TextValue 'E41' '  +0.11%  '

Set-TextValue 'D42' '103.30'
Set-TextValue 'E42' '  -5.83%  '

Set-TextValue 'D43' '0.01483'
Set-TextValue 'E43' '  -14.31%  '

Set-TextValue 'D44' '0.7834'
Set-TextValue 'E44' '  -12.95%  '

Set-TextValue 'B45' 'FraxShare'
Set-TextValue 'C45' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D45' '5.132'
Set-TextValue 'E45' '  -13.76%  '

Set-TextValue 'B46' 'TheSandbox'
Set-TextValue 'C46' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D46' '0.3788'
Set-TextValue 'E46' '  -21.05%  '

Set-TextValue 'B47' 'Cronos'
Set-TextValue 'C47' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D47' '0.05098'
Set-TextValue 'E47' '  -12.27%  '

Set-TextValue 'B48' 'Aptos'
Set-TextValue 'C48' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D48' '5.980'
Set-TextValue 'E48' '  -22.21%  '

Set-TextValue 'B49' 'Algorand'
Set-TextValue 'C49' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D49' '0.1061'
Set-TextValue 'E49' '  -14.36%  '

Set-TextValue 'B50' 'Elrond'
Set-TextValue 'C50' 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue 'D50' '30.12'
Set-TextValue 'E50' '  -13.36%  '

Set-TextValue 'D51' '52.34'
Set-TextValue 'E51' '  -12.67%  '
